$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "tblStudyFields1"

$data = @(
  @('Order', 'Field Name', 'Field Category', 'Field Type', 'Description'),
  @(47, 'MPManagementTool', 'Alternatives', 'List', 'Management alternatives evaluated in the study'),
  @(48, 'MPAlternativesEvaluated', 'Alternatives', 'List', 'Types of alternatives evaluated'),
  @(40, 'ObjName', 'Objectives', 'Description', 'Text description of the objective'),
  @(41, 'ObjCategory', 'Objectives', 'List', 'The objective category (conservation, yield, economic, social)'),
  @(42, 'ObjDescription', 'Objectives', 'Description', 'Description of the objective'),
  @(43, 'ObjDirection', 'Objectives', 'List', 'The desired state of the objective'),
  @(44, 'ObjType', 'Objectives', 'List', 'The type of objective.  E.g., strategic, process, fundamental, or means'),
  @(45, 'ObjScale', 'Objectives', 'Description', 'The scale on which the objective is measures (natural, proxy, or constructed)'),
  @(46, 'ObjMetric', 'Objectives', 'Description', 'The units used to measure the objective'),
  @(19, 'ResultsAdopted', 'Decision Process', 'Yes/No', 'Did the MSE influence subsequent management?'),
  @(20, 'Decision', 'Result', 'Description', 'If documented, the management procedure that was selected for implementation'),
  @(33, 'Drivers', 'System', 'List', 'System drivers included in the operating model, e.g., climate change, environmental conditions, predation, species interactions, etc.'),
  @(6, 'Latitude', 'Documentation', 'Lat.', 'Center point of the geographic are the MSE applies to (latitude)'),
  @(7, 'Longitude', 'Documentation', 'Long.', 'Center point of the geographic are the MSE applies to (longitude)'),
  @(2, 'Author', 'Documentation', 'Description', 'Names of the authors of the study'),
  @(3, 'YearPub', 'Documentation', 'Description', 'The year of publication'),
  @(8, 'Species', 'System', 'Description', 'Target species'),
  @(5, 'Location', 'System', 'Description', 'Geographic area the fishery occurs in'),
  @(4, 'System', 'System', 'Description', 'A short description of the study system.  For example, the species and it''s location'),
  @(9, 'FullCitation', 'Documentation', 'Description', 'The full citation for the study'),
  @(14, 'ProcessExplicit', 'Decision Analysis', 'Yes/No', 'Was the decision process methodology documented clearly?'),
  @(15, 'RoleSpecification', 'Decision Process', 'Yes/No', 'Were roles assigned and documented?'),
  @(16, 'OpenMeetings', 'Decision Process', 'Yes/No', 'Were open meetings held?'),
  @(23, 'ProblemDefinitionExplicit', 'Decision Analysis', 'Yes/No', 'Was a problem defintion completed and documented?'),
  @(24, 'ProblemDefinition', 'Decision Analysis', 'Description', 'A problem definition taken from the documentation (Reader interpretation)'),
  @(25, 'ObjectivesExplicit', 'Decision Analysis', 'Yes/No', 'Were objectives elicited and documented clearly?'),
  @(28, 'ObjElicitationMethod', 'Decision Analysis', 'Description', 'If documented, how were objectives elicited?'),
  @(34, 'TradeOffsExplicit', 'Decision Analysis', 'Yes/No', 'Was a tradeoff analysis conducted and documented?'),
  @(35, 'TradeOffMethod_Exp', 'Decision Analysis', 'Description', 'If explicitly documented, what form of tradeoff analysis occured?'),
  @(18, 'DecisionExplicit', 'Decision Analysis', 'Yes/No', 'Was the decision of the process documented?'),
  @(37, 'Comments', 'Comments', 'Description', 'Additional notes and comments about the study'),
  @(21, 'Leader', 'Decision Process', 'List', 'What organization initiated and directed the MSE?'),
  @(22, 'Participants', 'Decision Process', 'List', 'Who participated in the MSE process?'),
  @(26, 'ObjElicitationSource_Exp', 'Decision Process', 'List', 'If explicitly documented, the groups from which objectives were elicited'),
  @(30, 'ProcedureElicitation_Exp', 'Decision Process', 'List', 'If explicitly documented, the groups from which alternative management procedures were elicited'),
  @(32, 'ConsequencePrediction', 'Decision Process', 'Description', 'How were consequences predicted?'),
  @(17, 'OptimalAltExplicit', 'Decision Analysis', 'Yes/No', 'Was the best management procedure, aka optimal alternative, documented?'),
  @(29, 'AlternativesExplicit', 'Decision Analysis', 'Yes/No', 'Were alternatives elicited and documented clearly?'),
  @(36, 'TradeOffMethod_Sub', 'Decision Analysis', 'Description', 'If not explicitly documented, what form of tradeoff analysis seems to have occured?'),
  @(27, 'ObjElicitationSource_Sub', 'Decision Process', 'List', 'If not explicitly documented, the groups from which objectives were seemingly elicited'),
  @(31, 'ProcedureElicitation_Sub', 'Decision Process', 'List', 'If not explicitly documented, the groups from which alternative management procedures were seemingly elicited'),
  @(1, 'DOI', 'Documentation', 'ID', 'Digital Object Identifier (DOI)'),
  @(10, 'AllAuthors', 'Documentation', 'List', 'A list of all of the authors of the publication, in last, first, initial.; format'),
  @(11, 'Title', 'Documentation', 'Description', 'Title of the publication'),
  @(12, 'Journal', 'Documentation', 'Description', 'Journal in which the MSE study was published'),
  @(13, 'poc', 'Documentation', 'Description', 'Contact information (i.e., email address) for the publication''s point of contact'),
  @(38, 'RandomSample', 'Meta', 'Yes/No', 'Denotes whether the MSE was randomly sampled for review for inclusion in the publication, or not.'),
  @(39, 'UseInPublication', 'Meta', 'Yes/No', 'Denotes whether the MSE was reviewed for inclusion in the publication, or not.')
)

for ($r = 0; $r -lt $data.Count; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}

$wb.Names.Add('tblStudyFields', "='tblStudyFields1'!`$A`$1:`$E`$49") | Out-Null

$ws2.Range("C5").Select() | Out-Null
$ws2.Activate() | Out-Null
